$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A41").Value = "2025/12/04 09:00"
$ws.Range("B41").Value = "-"
$ws.Range("C41").Value = "-"
$ws.Range("D41").Value = "-"
$ws.Range("E41").Value = "-"
$ws.Range("F41").Value = "-"
$ws.Range("G41").Value = "-"
